$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header row): text labels, using the same bordered / bold /
# centered style that the header row already carries (style index 1) -------

# Stash the existing header formatting on a scratch cell (H1) before any
# header cell gets overwritten, so it can be re-applied afterwards without
# minting a brand-new style for every header cell (incl. the two new ones,
# D1/E1).
$ws.Range("C1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# Build each label as literal text in a scratch cell (G1) - the leading "'"
# is the normal Excel "treat as text" entry, same as typing it by hand - then
# copy value+type only (xlPasteAll) into the real destination.
$ws.Range("G1").Value = "'-3"
$ws.Range("G1").Copy()
$ws.Range("A1").PasteSpecial(-4104)

$ws.Range("G1").Value = "'1"
$ws.Range("G1").Copy()
$ws.Range("B1").PasteSpecial(-4104)

$ws.Range("G1").Value = "'0"
$ws.Range("G1").Copy()
$ws.Range("C1").PasteSpecial(-4104)

$ws.Range("G1").Value = "'0.1"
$ws.Range("G1").Copy()
$ws.Range("D1").PasteSpecial(-4104)

$ws.Range("G1").Value = "'0.2"
$ws.Range("G1").Copy()
$ws.Range("E1").PasteSpecial(-4104)

# Re-apply the stashed header formatting to all five header cells at once.
$ws.Range("H1").Copy()
$ws.Range("A1:E1").PasteSpecial(-4122)

# Drop the scratch cells - they must not show up in the final sheet.
$ws.Range("G1").Clear()
$ws.Range("H1").Clear()

# --- Row 2: numeric values (existing row, values updated) ------------------
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = -3
$ws.Range("C2").Value = -0.8
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0

# --- Row 3: numeric values (new row) ----------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = -10
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 1.95
$ws.Range("E3").Value = 0

# --- Row 4: numeric values (new row) ----------------------------------------
$ws.Range("A4").Value = -4
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = -2
$ws.Range("D4").Value = 1.5
$ws.Range("E4").Value = 0.45
